$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.95913233839191
$ws.Range("C2").Value = 7.017072097323165
$ws.Range("D2").Value = 4.332044269402864
$ws.Range("E2").Value = 10.97399696776474
$ws.Range("F2").Value = 62.91472073251706
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.79928908514583
$ws.Range("K2").Value = 14.64954740462116
$ws.Range("M2").Value = 17.18926605683532

$ws.Range("B3").Value = 14.91463669885406
$ws.Range("C3").Value = 7.031843805192101
$ws.Range("D3").Value = 4.328982487646406
$ws.Range("E3").Value = 11.02401195914878
$ws.Range("F3").Value = 61.96255502277742
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.79514011089733
$ws.Range("K3").Value = 14.65813803374158
$ws.Range("M3").Value = 17.24774409712304

$ws.Range("B4").Value = 14.89297276187163
$ws.Range("C4").Value = 7.044718253212088
$ws.Range("D4").Value = 4.329218401015022
$ws.Range("E4").Value = 11.05697670231988
$ws.Range("F4").Value = 61.3732326644574
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.79347060735602
$ws.Range("K4").Value = 14.66854854221333
$ws.Range("M4").Value = 17.28810834004006

$ws.Range("B5").Value = 14.88557416781955
$ws.Range("C5").Value = 7.050915421565786
$ws.Range("D5").Value = 4.329843502724875
$ws.Range("E5").Value = 11.07097758293293
$ws.Range("F5").Value = 61.13209582510508
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.79301135629655
$ws.Range("K5").Value = 14.67408013766446
$ws.Range("M5").Value = 17.30567587184844

$ws.Range("B6").Value = 14.88443215478903
$ws.Range("C6").Value = 7.052001686969665
$ws.Range("D6").Value = 4.329979144673646
$ws.Range("E6").Value = 11.07333670371981
$ws.Range("F6").Value = 61.09200161527794
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.79294845908231
$ws.Range("K6").Value = 14.67507643789968
$ws.Range("M6").Value = 17.30866044715385

$ws.Range("B7").Value = 14.89286718532389
$ws.Range("C7").Value = 7.044797989513117
$ws.Range("D7").Value = 4.329224694202368
$ws.Range("E7").Value = 11.05716322492073
$ws.Range("F7").Value = 61.36998432966059
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.79346351819122
$ws.Range("K7").Value = 14.66861792690259
$ws.Range("M7").Value = 17.28834073529964

$ws.Range("B8").Value = 14.94262119339596
$ws.Range("C8").Value = 7.021372430197369
$ws.Range("D8").Value = 4.330547781455425
$ws.Range("E8").Value = 10.9907742537397
$ws.Range("F8").Value = 62.58750389088479
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.79767648985795
$ws.Range("K8").Value = 14.6514424606184
$ws.Range("M8").Value = 17.20850255171219

$ws.Range("B9").Value = 15.08465530129254
$ws.Range("C9").Value = 7.005844693441661
$ws.Range("D9").Value = 4.350033520420967
$ws.Range("E9").Value = 10.87846527384856
$ws.Range("F9").Value = 64.92940524728276
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.81289058735868
$ws.Range("K9").Value = 14.65856863088584
$ws.Range("M9").Value = 17.08741666968587

$ws.Range("B10").Value = 15.2154216527304
$ws.Range("C10").Value = 7.013211513086044
$ws.Range("D10").Value = 4.374735681076705
$ws.Range("E10").Value = 10.80683127301626
$ws.Range("F10").Value = 66.61115317576689
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.82828760090488
$ws.Range("K10").Value = 14.68869659533474
$ws.Range("M10").Value = 17.02021264222937

$ws.Range("B11").Value = 15.28045771360501
$ws.Range("C11").Value = 7.020668535059809
$ws.Range("D11").Value = 4.388231212685513
$ws.Range("E11").Value = 10.77660156328593
$ws.Range("F11").Value = 67.36555461684833
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.83620324420451
$ws.Range("K11").Value = 14.70778971111341
$ws.Range("M11").Value = 16.99438911633913

$ws.Range("B12").Value = 15.30586463591695
$ws.Range("C12").Value = 7.024083403255556
$ws.Range("D12").Value = 4.393665349359692
$ws.Range("E12").Value = 10.76549305600936
$ws.Range("F12").Value = 67.64951859093711
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.83933119295665
$ws.Range("K12").Value = 14.71579127784135
$ws.Range("M12").Value = 16.98529491896504

$ws.Range("B13").Value = 15.30035846457877
$ws.Range("C13").Value = 7.023321664539921
$ws.Range("D13").Value = 4.392480640949487
$ws.Range("E13").Value = 10.76787040371465
$ws.Range("F13").Value = 67.58844081170034
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.83865174085592
$ws.Range("K13").Value = 14.71403375253758
$ws.Range("M13").Value = 16.98722304109123

$ws.Range("B14").Value = 15.28253245452269
$ws.Range("C14").Value = 7.02093763296531
$ws.Range("D14").Value = 4.38867180449172
$ws.Range("E14").Value = 10.7756808695923
$ws.Range("F14").Value = 67.38895182275677
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.83645797132671
$ws.Range("K14").Value = 14.70843256730385
$ws.Range("M14").Value = 16.99362720060026

$ws.Range("B15").Value = 15.2717143919854
$ws.Range("C15").Value = 7.019554317757096
$ws.Range("D15").Value = 4.386380888333656
$ws.Range("E15").Value = 10.78050912823074
$ws.Range("F15").Value = 67.26653088548144
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.83513119595061
$ws.Range("K15").Value = 14.70510202406604
$ws.Range("M15").Value = 16.99763914199302

$ws.Range("B16").Value = 15.21128149402546
$ws.Range("C16").Value = 7.012806879251306
$ws.Range("D16").Value = 4.373899086362892
$ws.Range("E16").Value = 10.80885427462618
$ws.Range("F16").Value = 66.5616215333402
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.82778858587443
$ws.Range("K16").Value = 14.6875569697854
$ws.Range("M16").Value = 17.02199597449393

$ws.Range("B17").Value = 15.17561641167643
$ws.Range("C17").Value = 7.009719954755139
$ws.Range("D17").Value = 4.366819813076643
$ws.Range("E17").Value = 10.8268466992878
$ws.Range("F17").Value = 66.12632684969367
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.82351717019872
$ws.Range("K17").Value = 14.67817154339732
$ws.Range("M17").Value = 17.03815554583895

$ws.Range("B18").Value = 15.15562641345305
$ws.Range("C18").Value = 7.008330928580826
$ws.Range("D18").Value = 4.362960688275034
$ws.Range("E18").Value = 10.83741731868932
$ws.Range("F18").Value = 65.87496639413203
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.82114619572299
$ws.Range("K18").Value = 14.67328068228145
$ws.Range("M18").Value = 17.04789690101494

$ws.Range("B19").Value = 15.14894862210561
$ws.Range("C19").Value = 7.007926970480566
$ws.Range("D19").Value = 4.361690604066507
$ws.Range("E19").Value = 10.841034455544
$ws.Range("F19").Value = 65.78969562136862
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.82035818344391
$ws.Range("K19").Value = 14.67171194932359
$ws.Range("M19").Value = 17.05127183437351

$ws.Range("B20").Value = 15.17935896008617
$ws.Range("C20").Value = 7.010008554676035
$ws.Range("D20").Value = 4.367551406989032
$ws.Range("E20").Value = 10.82490841491745
$ws.Range("F20").Value = 66.17276856505609
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.82396299069702
$ws.Range("K20").Value = 14.67911814628829
$ws.Range("M20").Value = 17.03638907611009

$ws.Range("B21").Value = 15.28774739776668
$ws.Range("C21").Value = 7.02162184014426
$ws.Range("D21").Value = 4.389781780227185
$ws.Range("E21").Value = 10.77337755143866
$ws.Range("F21").Value = 67.4475944402047
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.83709879885744
$ws.Range("K21").Value = 14.71005686691027
$ws.Range("M21").Value = 16.99172754961661

$ws.Range("B22").Value = 15.3631172553642
$ws.Range("C22").Value = 7.032656428260353
$ws.Range("D22").Value = 4.406196171289725
$ws.Range("E22").Value = 10.74167406913028
$ws.Range("F22").Value = 68.27071173082427
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.84644394948249
$ws.Range("K22").Value = 14.73477123625361
$ws.Range("M22").Value = 16.96652983068754

$ws.Range("B23").Value = 15.32248294572302
$ws.Range("C23").Value = 7.026451952569928
$ws.Range("D23").Value = 4.397263531778823
$ws.Range("E23").Value = 10.75841413685057
$ws.Range("F23").Value = 67.83237690245794
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.84138692654826
$ws.Range("K23").Value = 14.72117087259008
$ws.Range("M23").Value = 16.97961256764845

$ws.Range("B24").Value = 15.17766535043624
$ws.Range("C24").Value = 7.009876877325054
$ws.Range("D24").Value = 4.367219996569795
$ws.Range("E24").Value = 10.82578400763698
$ws.Range("F24").Value = 66.15177568946852
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.82376117132661
$ws.Range("K24").Value = 14.67868861401263
$ws.Range("M24").Value = 17.03718629237492

$ws.Range("B25").Value = 15.04153471362989
$ws.Range("C25").Value = 7.006761384674559
$ws.Range("D25").Value = 4.342941603911302
$ws.Range("E25").Value = 10.90693588967875
$ws.Range("F25").Value = 64.30190186393243
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.80803222904908
$ws.Range("K25").Value = 14.65226651243993
$ws.Range("M25").Value = 17.11636152672905
